$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) The existing "总计" sheet (3rd sheet) becomes the new "2022-Q1" sheet:
#    its old summary-of-quarters content is cleared and replaced with the
#    2022-Q1 fund holdings detail (same shape as the 2021-Q3 / 2021-Q4 sheets).
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(3)
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# Style reference cells (style index 2 on this workbook) to stamp onto the
# new header row + index column via copy/paste-special of formats only.
$styleSrc = $wb.Worksheets.Item(2).Range("B1")
$idxStyleSrc = $wb.Worksheets.Item(2).Range("A2")

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $q1.Cells.Item(1, 2 + $i)
    $cell.Value = $headers[$i]
}
$q1.Range("B1:H1").NumberFormat = "General"
$styleSrc.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$rows = @(
    @("519702", "交银趋势优先混合", "108.29", "71.40", "2.14", "2.3174", 7),
    @("010963", "信达澳银周期动力混合", "30.52", "89.82", "5.26", "1.6054", 3),
    @("010363", "信达澳银匠心臻选两年持有期混合", "50.40", "92.98", "2.11", "1.0634", 5),
    @("011160", "富国质量成长6个月持有期混合A", "3.80", "91.55", "2.95", "0.1121", 3),
    @("233009", "大摩多因子精选策略混合", "6.77", "89.73", "1.08", "0.0731", 7),
    @("001097", "华泰柏瑞积极优选股票", "1.25", "86.06", "3.00", "0.0375", 7),
    @("011161", "富国质量成长6个月持有期混合C", "0.12", "91.55", "2.95", "0.0035", 3)
)

# Columns whose values must round-trip as plain text (not auto-coerced to
# numbers, which would both change the stored type AND strip leading zeros
# from fund codes like "010963"). NumberFormat="@" forces text entry; then
# resetting the Style back to "Normal" drops the formatting so the cell ends
# up with no explicit style index, matching the target cells exactly.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowNum = 2 + $r
    $data = $rows[$r]

    $idxCell = $q1.Cells.Item($rowNum, 1)
    $idxCell.Value = $r
    $idxStyleSrc.Copy()
    $idxCell.PasteSpecial(-4122)

    Set-TextValue $q1.Cells.Item($rowNum, 2) $data[0]
    $q1.Cells.Item($rowNum, 3).Value = $data[1]

    Set-TextValue $q1.Cells.Item($rowNum, 4) $data[2]
    Set-TextValue $q1.Cells.Item($rowNum, 5) $data[3]
    Set-TextValue $q1.Cells.Item($rowNum, 6) $data[4]
    Set-TextValue $q1.Cells.Item($rowNum, 7) $data[5]

    $q1.Cells.Item($rowNum, 8).Value = $data[6]
}

# ---------------------------------------------------------------------------
# 2) A brand-new "总计" sheet is appended at the end, carrying forward the
#    old per-quarter summary rows plus a new 2022-Q1 entry on top.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Match the page-margin convention used by the rest of this workbook's
# sheets (0.75"/0.75"/1"/1"/0.5"/0.5" -- points = inches * 72) rather than
# the host's brand-new-sheet defaults (0.7"/0.7"/0.75"/0.75"/0.3"/0.3").
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"
$total.Range("B1:D1").NumberFormat = "General"
$styleSrc.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$totalRows = @(
    @("2022-Q1", 7, 5.21),
    @("2021-Q4", 3, 4.33),
    @("2021-Q3", 1, 2.44)
)

for ($r = 0; $r -lt $totalRows.Count; $r++) {
    $rowNum = 2 + $r
    $data = $totalRows[$r]

    $idxCell = $total.Cells.Item($rowNum, 1)
    $idxCell.Value = $r
    $idxStyleSrc.Copy()
    $idxCell.PasteSpecial(-4122)

    $total.Cells.Item($rowNum, 2).Value = $data[0]
    $total.Cells.Item($rowNum, 3).Value = $data[1]
    $total.Cells.Item($rowNum, 4).Value = $data[2]
}
